$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 79.541664
$ws.Range("I9").Value = 85.388885
$ws.Range("K9").Value = 85.388885
$ws.Range("M9").Value = 83.611115
$ws.Range("H15").Value = 514.8333
$ws.Range("I15").Value = 514.8333
$ws.Range("K15").Value = 1544.4999
$ws.Range("M15").Value = -1375.4999
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""
$ws.Range("H98").Value = 8981.200000000001
$ws.Range("I98").Value = 5114.2
$ws.Range("J98").Value = 12848.2
$ws.Range("K98").Value = 5114.2
$ws.Range("L98").Value = 12848.2
$ws.Range("M98").Value = -3616.2
$ws.Range("N98").Value = -15844.2
$ws.Range("H113").Value = 2764.1428
$ws.Range("I113").Value = 2119.8
$ws.Range("J113").Value = 4375
$ws.Range("K113").Value = 2119.8
$ws.Range("L113").Value = 4375
$ws.Range("M113").Value = 1134.2
$ws.Range("N113").Value = -10883
$ws.Range("H122").Value = 8981.200000000001
$ws.Range("I122").Value = 5114.2
$ws.Range("J122").Value = 12848.2
$ws.Range("K122").Value = 15342.6
$ws.Range("L122").Value = 38544.60000000001
$ws.Range("M122").Value = -12892.6
$ws.Range("N122").Value = -43444.60000000001
$ws.Range("H138").Value = 4765242.5
$ws.Range("J138").Value = 3690.2
$ws.Range("L138").Value = 11070.6
$ws.Range("N138").Value = -21350.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 536
$ws.Range("I4").Value = 191
$ws.Range("K4").Value = 191
$ws.Range("M4").Value = -75
$ws.Range("H63").Value = 4406.8335
$ws.Range("I63").Value = 4548.8
$ws.Range("J63").Value = 3697
$ws.Range("K63").Value = 4548.8
$ws.Range("L63").Value = 3697
$ws.Range("M63").Value = -3862.8
$ws.Range("N63").Value = -5069
$ws.Range("H66").Value = 4406.8335
$ws.Range("I66").Value = 4548.8
$ws.Range("J66").Value = 3697
$ws.Range("K66").Value = 22744
$ws.Range("L66").Value = 18485
$ws.Range("M66").Value = -19312
$ws.Range("N66").Value = -25349
$ws.Range("H132").Value = 1584.7693
$ws.Range("I132").Value = 1080.2
$ws.Range("K132").Value = 3240.6
$ws.Range("M132").Value = -710.6000000000004
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1734.1
$ws.Range("J134").Value = 2749.5
$ws.Range("L134").Value = 8248.5
$ws.Range("N134").Value = -13318.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 46.5
$ws.Range("I7").Value = 46.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 46.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 66.5
$ws.Range("N7").Value = ""
$ws.Range("H17").Value = 10000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("H22").Value = 279.33334
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = ""
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("H32").Value = 4049.6
$ws.Range("I32").Value = 2875.5
$ws.Range("K32").Value = 2875.5
$ws.Range("M32").Value = -2559.5
$ws.Range("H41").Value = 8000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H50").Value = 10000
$ws.Range("I50").Value = 10000
$ws.Range("K50").Value = 10000
$ws.Range("M50").Value = -9375
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 6000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5264
$ws.Range("H58").Value = 2942.625
$ws.Range("I58").Value = 3077.5715
$ws.Range("J58").Value = 1998
$ws.Range("K58").Value = 3077.5715
$ws.Range("L58").Value = 1998
$ws.Range("M58").Value = -2874.5715
$ws.Range("N58").Value = -2404
$ws.Range("H60").Value = 27200
$ws.Range("I60").Value = 27200
$ws.Range("K60").Value = 27200
$ws.Range("M60").Value = -26689
$ws.Range("H61").Value = 6000
$ws.Range("I61").Value = 6000
$ws.Range("K61").Value = 6000
$ws.Range("M61").Value = -5652
$ws.Range("H107").Value = 762.8182
$ws.Range("I107").Value = 539.1
$ws.Range("K107").Value = 539.1
$ws.Range("M107").Value = 1380.9
$ws.Range("H122").Value = 2672.3333
$ws.Range("I122").Value = 2759
$ws.Range("K122").Value = 8277
$ws.Range("M122").Value = -5827
$ws.Range("H136").Value = 2942.625
$ws.Range("I136").Value = 3077.5715
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 9232.7145
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -6682.7145
$ws.Range("N136").Value = -11094
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 500
$ws.Range("J23").Value = 500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1970
$ws.Range("H34").Value = 698.5
$ws.Range("I34").Value = 530.8333
$ws.Range("K34").Value = 1592.4999
$ws.Range("M34").Value = -1508.4999
$ws.Range("H40").Value = 191.33333
$ws.Range("J40").Value = 161.5
$ws.Range("L40").Value = 646
$ws.Range("N40").Value = -784
$ws.Range("H48").Value = 404
$ws.Range("J48").Value = 404
$ws.Range("L48").Value = 1212
$ws.Range("N48").Value = -1712
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = ""
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = ""
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = ""
$ws.Range("H104").Value = 4592.25
$ws.Range("I104").Value = 400
$ws.Range("K104").Value = 1200
$ws.Range("M104").Value = 1421
$ws.Range("H131").Value = 14999.5
$ws.Range("I131").Value = 14999.5
$ws.Range("K131").Value = 44998.5
$ws.Range("M131").Value = -39958.5
$ws.Range("H137").Value = 1916.5
$ws.Range("J137").Value = 2033
$ws.Range("L137").Value = 6099
$ws.Range("N137").Value = -16299
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 24999
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = ""
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1886.6
$ws.Range("I61").Value = 1644.6666
$ws.Range("J61").Value = 2249.5
$ws.Range("K61").Value = 1644.6666
$ws.Range("L61").Value = 2249.5
$ws.Range("M61").Value = -1442.6666
$ws.Range("N61").Value = -2653.5
$ws.Range("H93").Value = 1460.2667
$ws.Range("I93").Value = 1300.3846
$ws.Range("K93").Value = 1300.3846
$ws.Range("M93").Value = -52.38460000000009
$ws.Range("H113").Value = 1886.6
$ws.Range("I113").Value = 1644.6666
$ws.Range("J113").Value = 2249.5
$ws.Range("K113").Value = 1644.6666
$ws.Range("L113").Value = 2249.5
$ws.Range("M113").Value = 525.3334
$ws.Range("N113").Value = -6589.5
$ws.Range("H132").Value = 3688.4167
$ws.Range("I132").Value = 2526.1
$ws.Range("K132").Value = 7578.299999999999
$ws.Range("M132").Value = -5048.299999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 29999.5
$ws.Range("H43").Value = 17250
$ws.Range("J43").Value = 17250
$ws.Range("L43").Value = 17250
$ws.Range("N43").Value = -17548
$ws.Range("H52").Value = 19666
$ws.Range("I52").Value = 19666
$ws.Range("K52").Value = 19666
$ws.Range("M52").Value = -19440
$ws.Range("H62").Value = 135833
$ws.Range("I62").Value = 7499
$ws.Range("J62").Value = 200000
$ws.Range("K62").Value = 7499
$ws.Range("L62").Value = 200000
$ws.Range("M62").Value = -6875
$ws.Range("N62").Value = -201248
$ws.Range("H65").Value = 135833
$ws.Range("I65").Value = 7499
$ws.Range("J65").Value = 200000
$ws.Range("K65").Value = 37495
$ws.Range("L65").Value = 1000000
$ws.Range("M65").Value = -34375
$ws.Range("N65").Value = -1006240
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = ""
$ws.Range("N132").Value = ""
